# Commit: "Batch 42 Dec 10 2018"
# Duplicate the "Sept 7 2017" source-calibration-check sheet, place the
# copy right after it, rename it to "Dec 10 2018", and fill in the new
# calibration data for CCSEO Batch No. 42.

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Sept 7 2017")

# Duplicate the sheet, inserting the new copy immediately after the source.
$src.Copy($null, $src)
$ws = $wb.Worksheets.Item($src.Index + 1)
$ws.Name = "Dec 10 2018"

# Batch label.
$ws.Range("A3").Value = "CCSEO Batch No.: 42"

# Manufacturer's activity + calibration date/time.
$ws.Range("B4").Value = 43431.072916666664

# Measured current (nA).
$ws.Range("B6").Value = 383.8

# Local calibration date/time.
$ws.Range("C8").Value = 43444.620138888888

# Temperature / pressure readings.
$ws.Range("H8").Value = 21.8
$ws.Range("H9").Value = 762.2
$ws.Range("C10").Value = 1.0008999999999999

$ws.Range("C11").Value = 462.1
$ws.Range("C12").Value = -80.367999999999995

# Console activity expected value.
$ws.Range("D19").Value = 37092

# Make sure this is the print area for the new sheet (named range
# _xlnm.Print_Area, local to this sheet).
$ws.PageSetup.PrintArea = 'A1:H22'

# Keep the new sheet as the active / selected tab.
$ws.Select()
